$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename TestBean to JavaBean to avoid interpolation with JUnit tests in surefire plugin
$ws.Range("B3").Value = "Method String print(JavaBean bean)"
$ws.Range("F4").Value = "Data JavaBean beans"

# Move active selection to F5 (as reflected in the saved workbook)
$ws.Range("F5").Select()
